$d = $word.ActiveDocument

$pairs = @(
    @("59÷7=8, 3", "48÷2=24, 0"),
    @("12÷2=6, 0", "48÷2=24, 0"),
    @("14÷8=1, 6", "75÷2=37, 1"),
    @("40÷8=5, 0", "33÷2=16, 1"),
    @("87÷7=12, 3", "89÷6=14, 5"),
    @("67÷2=33, 1", "48÷2=24, 0"),
    @("97÷4=24, 1", "32÷4=8, 0"),
    @("29÷3=9, 2", "87÷6=14, 3"),
    @("64÷8=8, 0", "77÷4=19, 1"),
    @("79÷7=11, 2", "37÷4=9, 1"),
    @("46÷7=6, 4", "58÷7=8, 2"),
    @("75÷7=10, 5", "16÷3=5, 1"),
    @("80÷3=26, 2", "86÷8=10, 6"),
    @("28÷4=7, 0", "55÷7=7, 6"),
    @("66÷2=33, 0", "25÷8=3, 1"),
    @("95÷4=23, 3", "30÷2=15, 0"),
    @("65÷5=13, 0", "51÷9=5, 6"),
    @("81÷9=9, 0", "41÷8=5, 1"),
    @("28÷5=5, 3", "15÷4=3, 3"),
    @("53÷3=17, 2", "46÷5=9, 1"),
    @("85÷6=14, 1", "65÷4=16, 1"),
    @("63÷8=7, 7", "94÷6=15, 4"),
    @("60÷9=6, 6", "57÷9=6, 3"),
    @("51÷2=25, 1", "95÷3=31, 2"),
    @("84÷6=14, 0", "28÷6=4, 4")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

$d.Save()
